# Update "想去人数" (want-to-go count, column F) figures on the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to match the
# latest scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" updates
$ws1.Range("F2").Value = 5529
$ws1.Range("F4").Value = 7640
$ws1.Range("F6").Value = 75
$ws1.Range("F8").Value = 611
$ws1.Range("F9").Value = 68
$ws1.Range("F10").Value = 37
$ws1.Range("F11").Value = 4402
$ws1.Range("F12").Value = 1787
$ws1.Range("F13").Value = 117
$ws1.Range("F14").Value = 117
$ws1.Range("F15").Value = 2979
$ws1.Range("F17").Value = 571
$ws1.Range("F18").Value = 221
$ws1.Range("F20").Value = 470
$ws1.Range("F21").Value = 478
$ws1.Range("F22").Value = 338
$ws1.Range("F24").Value = 1719
$ws1.Range("F25").Value = 1244
$ws1.Range("F27").Value = 1437
$ws1.Range("F29").Value = 593
$ws1.Range("F32").Value = 23
$ws1.Range("F33").Value = 20
$ws1.Range("F34").Value = 69
$ws1.Range("F37").Value = 3081
$ws1.Range("F38").Value = 716
$ws1.Range("F39").Value = 44
$ws1.Range("F40").Value = 139
$ws1.Range("F41").Value = 47
$ws1.Range("F42").Value = 797

# Sheet "全部类型" updates
$ws4.Range("F2").Value = 5529
$ws4.Range("F4").Value = 7640
$ws4.Range("F6").Value = 75
$ws4.Range("F8").Value = 611
$ws4.Range("F9").Value = 68
$ws4.Range("F10").Value = 37
$ws4.Range("F11").Value = 4402
$ws4.Range("F12").Value = 1787
$ws4.Range("F13").Value = 117
$ws4.Range("F14").Value = 117
$ws4.Range("F15").Value = 2979
$ws4.Range("F17").Value = 571
$ws4.Range("F18").Value = 221
$ws4.Range("F20").Value = 470
$ws4.Range("F21").Value = 478
$ws4.Range("F23").Value = 338
$ws4.Range("F25").Value = 1719
$ws4.Range("F26").Value = 1244
$ws4.Range("F28").Value = 1437
$ws4.Range("F30").Value = 593
$ws4.Range("F33").Value = 23
$ws4.Range("F34").Value = 20
$ws4.Range("F35").Value = 69
$ws4.Range("F38").Value = 3082
$ws4.Range("F40").Value = 716
$ws4.Range("F41").Value = 44
$ws4.Range("F42").Value = 139
$ws4.Range("F43").Value = 47
$ws4.Range("F44").Value = 797
